# Update column F ("dSF") values for the rows whose source data was
# re-pulled. Each entry maps an Excel row number to its new dSF value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 0
    4  = 1
    7  = 3
    10 = 1
    12 = 0
    13 = -3
    15 = 1
    16 = 1
    21 = 1
    22 = -3
    29 = -2
    33 = -2
    35 = 2
    38 = 3
    40 = -1
    41 = -1
    43 = 2
    50 = -2
    55 = -3
    56 = 6
    58 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
